# Update SWH (column B) values for rows 11-18 on Sheet1 — recomputed load
# stats (standard deviation / period) shifted the per-element SWH figures
# down by one bin, with a new value for element n16 (row 11).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B11").Value = 9.6
$ws.Range("B12").Value = 10.5
$ws.Range("B13").Value = 11.4
$ws.Range("B14").Value = 12.3
$ws.Range("B15").Value = 13.2
$ws.Range("B16").Value = 14.2
$ws.Range("B17").Value = 15.09
$ws.Range("B18").Value = 16

# Leave the cursor on the last touched cell, matching the author's
# final selection in the saved workbook.
$ws.Activate()
$ws.Range("D19").Select()
